$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("name","age","gender","phone","grade","parentName","parentPhone","studentPassword","studentUsername","parentUsername","parentPassword")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Student rows: name, age, gender, phone, grade, parentName, parentPhone, studentPassword, studentUsername, parentUsername, parentPassword
$rows = @(
    @("Rediet Erbelo", 15, "F", 251920864496, 9, "Betelhem Erbelo", 251920864496, "kafsxp", "redieterbelo4112", "betelhemerbelo7125", "opc1xc"),
    @("Kalkidan Erbelo", 10, "F", 251920864496, 5, "Eyu Erbelo", 251920864496, "v610q5", "kalkidanerbelo2678", "eyuerbelo1742", "zn631h"),
    @("Samuel Ayalew", 19, "M", 251931653440, 12, "Ayalew Bikago", 251931653440, "75jzd1", "samuelayalew3114", "ayalewbikago1793", "mf2yuz")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
